$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.182.33'
$ws.Range("E2").Value = '  +1.80%  '
$ws.Range("D3").Value = '2.024.65'
$ws.Range("E3").Value = '  +3.69%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''247.45'
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").Value = '''0.627'
$ws.Range("E6").Value = '  +1.88%  '
$ws.Range("D7").Value = '''60.08'
$ws.Range("E7").Value = '  -0.60%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '''0.393'
$ws.Range("E9").Value = '  +4.06%  '
$ws.Range("D10").Value = '''0.0806'
$ws.Range("E10").Value = '  +2.34%  '
$ws.Range("E11").Value = '  +1.96%  '
$ws.Range("D12").Value = '''15.14'
$ws.Range("E12").Value = '  +6.17%  '
$ws.Range("D13").Value = '2.323.51'
$ws.Range("E13").Value = '  +3.81%  '
$ws.Range("D14").Value = '''0.853'
$ws.Range("E14").Value = '  +3.41%  '
$ws.Range("D15").Value = '''22.03'
$ws.Range("E15").Value = '  +2.24%  '
$ws.Range("D16").Value = '''5.49'
$ws.Range("E16").Value = '  +4.81%  '
$ws.Range("D17").Value = '2.026.22'
$ws.Range("E17").Value = '  +3.99%  '
$ws.Range("D18").Value = '37.133.39'
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("D19").Value = '''70.44'
$ws.Range("E19").Value = '  +1.59%  '
$ws.Range("D20").Value = '0.0₃0863'
$ws.Range("E20").Value = '  +1.52%  '
$ws.Range("D21").Value = '''5.23'
$ws.Range("E21").Value = '  +3.19%  '
$ws.Range("D22").Value = '''230.52'
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").Value = '''2.59'
$ws.Range("E24").Value = '  +6.03%  '
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").Value = '''9.42'
$ws.Range("E26").Value = '  +3.02%  '
$ws.Range("D27").Value = '''162.94'
$ws.Range("E27").Value = '  +1.95%  '
$ws.Range("D28").Value = '''0.138'
$ws.Range("E28").Value = '  -4.17%  '
$ws.Range("D29").Value = '''19.73'
$ws.Range("E29").Value = '  +2.32%  '
$ws.Range("E30").Value = '  +4.36%  '
$ws.Range("D31").Value = '''0.121'
$ws.Range("E31").Value = '  +1.69%  '
$ws.Range("D32").Value = '''4.80'
$ws.Range("E32").Value = '  +1.53%  '
$ws.Range("D33").Value = '''0.0672'
$ws.Range("E33").Value = '  +9.81%  '
$ws.Range("E34").Value = '  +9.41%  '
$ws.Range("D35").Value = '''4.47'
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("D36").Value = '''3.64'
$ws.Range("E36").Value = '  +4.92%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  +2.35%  '
$ws.Range("D39").Value = '''5.42'
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("D40").Value = '''3.04'
$ws.Range("E40").Value = '  +4.06%  '
$ws.Range("D41").Value = '''0.0980'
$ws.Range("E41").Value = '  +1.77%  '
$ws.Range("D42").Value = '''17.05'
$ws.Range("E42").Value = '  +8.41%  '
$ws.Range("D43").Value = '''1.18'
$ws.Range("E43").Value = '  +1.12%  '
$ws.Range("D44").Value = '''0.0213'
$ws.Range("E44").Value = '  +1.93%  '
$ws.Range("D45").Value = '''92.05'
$ws.Range("E45").Value = '  +4.06%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.378.56'
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '''1.06'
$ws.Range("E47").Value = '  +3.42%  '
$ws.Range("D48").Value = '''7.45'
$ws.Range("E48").Value = '  +4.70%  '
$ws.Range("E49").Value = '  +18.44%  '
$ws.Range("D50").Value = '''2.85'
$ws.Range("E50").Value = '  +0.78%  '
$ws.Range("D51").Value = '''46.06'
$ws.Range("E51").Value = '  +1.63%  '
